# Add data for 2023-11-15
# Updates YTD violent-crime counts across the citywide, by-neighborhood
# summary sheets and the individual neighborhood sheets affected by the
# new day of data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("Citywide Totals")
$ws.Range("F2").Value = 82
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 112
$ws.Range("E3").Value = 137
$ws.Range("H3").Value = 143
$ws.Range("J3").Value = 212
$ws.Range("E4").Value = 11
$ws.Range("J4").Value = 20
$ws.Range("D6").Value = 387
$ws.Range("E6").Value = 437
$ws.Range("F6").Value = 485
$ws.Range("G6").Value = 420
$ws.Range("I6").Value = 477
$ws.Range("J6").Value = 391
$ws.Range("D7").Value = 607
$ws.Range("E7").Value = 654
$ws.Range("F7").Value = 700
$ws.Range("G7").Value = 640
$ws.Range("H7").Value = 682
$ws.Range("I7").Value = 797
$ws.Range("J7").Value = 736

$ws = $wb.Worksheets("By Neighborhood")
$ws.Range("E2").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("E8").Value = 47
$ws.Range("F8").Value = 43
$ws.Range("G8").Value = 32
$ws.Range("J8").Value = 43
$ws.Range("F19").Value = 22
$ws.Range("G29").Value = 10
$ws.Range("E32").Value = 61
$ws.Range("E36").Value = 34
$ws.Range("E47").Value = 16
$ws.Range("E48").Value = 6
$ws.Range("D53").Value = 68
$ws.Range("E53").Value = 80
$ws.Range("H53").Value = 92
$ws.Range("I53").Value = 122
$ws.Range("J65").Value = 11
$ws.Range("J70").Value = 12
$ws.Range("J74").Value = 21
$ws.Range("D81").Value = 3
$ws.Range("E88").Value = 8
$ws.Range("D98").Value = 607
$ws.Range("E98").Value = 654
$ws.Range("F98").Value = 700
$ws.Range("G98").Value = 640
$ws.Range("H98").Value = 682
$ws.Range("I98").Value = 797
$ws.Range("J98").Value = 736

$ws = $wb.Worksheets("Austin")
$ws.Range("F2").Value = 8
$ws.Range("E5").Value = 37
$ws.Range("G5").Value = 23
$ws.Range("J5").Value = 23
$ws.Range("E6").Value = 47
$ws.Range("F6").Value = 43
$ws.Range("G6").Value = 32
$ws.Range("J6").Value = 43

$ws = $wb.Worksheets("Garfield Park")
$ws.Range("E6").Value = 48
$ws.Range("E7").Value = 61

$ws = $wb.Worksheets("Grand Crossing")
$ws.Range("E3").Value = 8
$ws.Range("E7").Value = 34

$ws = $wb.Worksheets("Washington Park")
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 8

$ws = $wb.Worksheets("Loop")
$ws.Range("H2").Value = 12
$ws.Range("I2").Value = 12
$ws.Range("H3").Value = 18
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 62
$ws.Range("I6").Value = 78
$ws.Range("D7").Value = 68
$ws.Range("E7").Value = 80
$ws.Range("H7").Value = 92
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets("South Deering")
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 3

$ws = $wb.Worksheets("North Lawndale")
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 11

$ws = $wb.Worksheets("River North")
$ws.Range("J3").Value = 5
$ws.Range("J6").Value = 21

$ws = $wb.Worksheets("Fuller Park")
$ws.Range("G5").Value = 8
$ws.Range("G6").Value = 10

$ws = $wb.Worksheets("Chatham")
$ws.Range("F5").Value = 15
$ws.Range("F6").Value = 22

$ws = $wb.Worksheets("Albany Park")
$ws.Range("E4").Value = 2
$ws.Range("E6").Value = 4

$ws = $wb.Worksheets("Lincoln Park")
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 6

$ws = $wb.Worksheets("Lake View")
$ws.Range("E5").Value = 10
$ws.Range("E6").Value = 16

$ws = $wb.Worksheets("Old Town")
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 12

$ws = $wb.Worksheets("Ashburn")
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
